# Fruta / hortaliza, semanal
# A new weekly price record is inserted at the top of the data table
# (row 219), pushing the existing rows 219:307 down to 220:308.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 219 - this shifts rows 219:307 down to 220:308
# and keeps the dimension / formatting in sync automatically.
$ws.Rows.Item(219).Insert()

# The freshly inserted row 219 is blank; the record that used to live at
# row 219 is now at row 220. Copy all of its columns down into the new
# row 219 first, then overwrite the handful of fields that differ for the
# new weekly observation (date, volume, prices, $/kg).
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(219, $col).Value = $ws.Cells.Item(220, $col).Value2
}

# New weekly entry values (everything else matches the carried-over row).
$ws.Cells.Item(219, 4).Value  = 44825   # D - Fecha
$ws.Cells.Item(219, 13).Value = 20      # M - Volumen
$ws.Cells.Item(219, 14).Value = 23000   # N - Precio minimo
$ws.Cells.Item(219, 15).Value = 23500   # O - Precio maximo
$ws.Cells.Item(219, 16).Value = 23250   # P - Precio promedio ponderado
$ws.Cells.Item(219, 19).Value = 1938    # S - Precio $/Kg
